# Refresh the cryptocurrency snapshot (Price column D, 1h change column E)
# with the latest scraped figures. Most of these look like plain numbers
# (e.g. "0.505", "19.74") but the source data stores the whole Price/Volume
# table as text, so any cell whose new value would otherwise be auto-
# converted to a Number by Excel is briefly switched to Text format, the
# text value is written, and the cell style is put back to "Normal" right
# away so no stray formatting is left behind on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.820.28'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.640.35'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '1.866.70'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = '1.635.74'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '25.857.50'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").Value = '  +2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.91%  '
$ws.Range("E24").Value = '  +5.74%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("D37").Value = '1.134.98'
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.808'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.10%  '
$ws.Range("D45").Value = '1.775.85'
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("E46").Value = '  +2.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.417'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.28%  '
